$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.943.54"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "3.438.99"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.17"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.53"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "3.439.33"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.74"
$ws.Range("E10").Value = "  +1.14%  "
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.403"
$ws.Range("E12").Value = "  +2.46%  "
$ws.Range("D13").Value = "4.030.23"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("E14").Value = "  +2.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.96"
$ws.Range("E15").Value = "  -1.62%  "
$ws.Range("D16").Value = "3.425.98"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "62.994.32"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.35"
$ws.Range("E19").Value = "  +1.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.36"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.19"
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "385.11"
$ws.Range("E22").Value = "  -2.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.559"
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.38"
$ws.Range("E24").Value = "  -1.31%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "3.585.02"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000114"
$ws.Range("E27").Value = "  -3.06%  "
$ws.Range("E28").Value = "  -5.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.57"
$ws.Range("E29").Value = "  -2.16%  "
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.08"
$ws.Range("E31").Value = "  -1.43%  "
$ws.Range("E32").Value = "  -2.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.26"
$ws.Range("E34").Value = "  -2.24%  "
$ws.Range("E35").Value = "  -8.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.28"
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.06"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "31.52"
$ws.Range("E39").Value = "  +2.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "169.20"
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("D41").Value = "3.477.44"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.787"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.33"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.19"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.34"
$ws.Range("E47").Value = "  -3.17%  "
$ws.Range("D48").Value = "2.563.09"
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.29"
$ws.Range("E49").Value = "  +5.62%  "
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.66"
$ws.Range("E51").Value = "  -4.08%  "
